$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$newValues = @(
    "36+43=",
    "83-65=",
    "17+56=",
    "1+83=",
    "28-16=",
    "91-68=",
    "84-9=",
    "92-68=",
    "57-50=",
    "65+0=",
    "11+36=",
    "15+48=",
    "93-13=",
    "78-49=",
    "18-12=",
    "99-71=",
    "58-40=",
    "35-28=",
    "38+37=",
    "78-25=",
    "18+27=",
    "13+85=",
    "54-18=",
    "30+44=",
    "27+52=",
    "13-5=",
    "17-2=",
    "16-13=",
    "49+12=",
    "87-16=",
    "92-81=",
    "79-34=",
    "47+4=",
    "28-8=",
    "20+54=",
    "96-79=",
    "60-47=",
    "21+21=",
    "90-59=",
    "57-30=",
    "58+40=",
    "88-44=",
    "51+14=",
    "10+49=",
    "86-76=",
    "62+4=",
    "81-76=",
    "3+16=",
    "38+17=",
    "11-2=",
    "74-54=",
    "29+11=",
    "19+57=",
    "64-3=",
    "89-18=",
    "43+8=",
    "77-58=",
    "77-30=",
    "71-5=",
    "61+6=",
    "25+42=",
    "84+6=",
    "11+39=",
    "28-14=",
    "61+10=",
    "13+62=",
    "50+47=",
    "45-15=",
    "57+7=",
    "81-37=",
    "82+9=",
    "49+8=",
    "0+72=",
    "68-47=",
    "5+20=",
    "3+96=",
    "77-18=",
    "42-21=",
    "82-15=",
    "60-27=",
    "9+49=",
    "32-10=",
    "26+60=",
    "88-8=",
    "34+10=",
    "23-5=",
    "11+33=",
    "53-11=",
    "4+70=",
    "80+19=",
    "98-19=",
    "63-27=",
    "66-48=",
    "77-31=",
    "72-53=",
    "78-74=",
    "52+45=",
    "83-18=",
    "50-28=",
    "87-75="
)

$rows = $tbl.Rows.Count
$cols = $tbl.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Host "Updated" $idx "cells"
